$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "D" column holds a shared formula that concatenates the hash in column C
# with a trailing comment built from column B ("C#&","&" // "&B#").
# Historically two sub-ranges (D17:D26/D34:D66 and D67:D95) were left on an
# older formula ("C#&",""), missing the " // <char>" comment suffix, while a
# handful of rows (D3:D16, D27:D33, D60:D65, D92:D95) already used the full
# formula. Bring the whole D3:D95 block back in sync with the full formula so
# every row gets its "// <char> <name>" comment again.
for ($r = 3; $r -le 95; $r++) {
    $ws.Range("D$r").Formula = "=C$r&"",""&"" // ""&B$r"
}

# Restore the selection to the top of the sheet (was scrolled down to A82 /
# D92:D95 after the previous edit).
$ws.Activate()
$ws.Range("D2:D95").Select()
